$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF ("Date") held values like "5-9-2007-08". Because of how the NBA
# stats site displayed dates, every row in this sheet (game date 2008-05-09)
# was off by a day / mis-formatted. Correct it to the proper ISO-style date
# text "2008-05-09" for all 30 data rows (BF2:BF31).
#
# NumberFormat is forced to Text ("@") before the write so the engine keeps
# the value as literal text instead of re-interpreting the ISO-looking
# string as a date serial number.
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF" + $r).Value = "2008-05-09"
}
